$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 11 (pushes "Description" and everything below down by one)
$ws.Rows.Item(11).Insert()

# The inserted row doesn't inherit the boxed/wrapped data-row style, so copy the
# formatting down from the row just below (row 12, a normal data row) onto the new row 11.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Populate the newly inserted row with the "Jurisdiction" property (empty value)
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Update Version value
$ws.Range("B3").Value = "0.1.1"

# Update Date value
$ws.Range("B8").Value = "2024-11-11T17:53:38-06:00"
